$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11) for every data row (2..451).
$ws.Range("C2:C451").Value = 45180
